# Implements the commit:
#   refactor(Docx + main): Generation of file body now is implemented into Docx
#
# Paragraph numbers below are 1-indexed positions in the *original* before.docx
# (Paragraphs.Item is 1-indexed). We perform same-count text replacements first,
# then the block replacement (which grows the paragraph count), and finally the
# single-paragraph insertion -- in that order so none of the fixed indices below
# are invalidated by an earlier step.

$d = $word.ActiveDocument

# --- 1) "use docx_rs::{...}" : drop the trailing ", SpecialIndentType" ---
$d.Paragraphs.Item(5).Range.Text = 'use docx_rs::{Docx, DocxError, IndentLevel, NumberingId, Paragraph, Run};'

# --- 2) call site: gen_file(path, &mut doc) -> doc.gen_body(path) ---
$d.Paragraphs.Item(51).Range.Text = '                    println!("{:?}", doc.gen_body(path));'

# --- 3) Replace the old free function `fn gen_file` (paragraphs 62-82, inclusive)
#     with the new `trait GenFile` / `impl GenFile for Docx` block (54 lines),
#     keeping the old body around, commented out, at the end.
#     The first 21 new lines reuse the 21 existing paragraphs in place;
#     the remaining 33 lines are appended as new paragraphs.
$genFileBlock = @(
    'trait GenFile {'
    '    fn gen_body(&mut self, input_path: PathBuf) -> Result<(), DocxError>;'
    '}'
    'impl GenFile for Docx {'
    '    fn gen_body(&mut self, input_path: PathBuf) -> Result<(), DocxError> {'
    '        *self = self.to_owned().add_paragraph('
    '            Paragraph::new()'
    '                .add_run('
    '                    Run::new()'
    '                        .add_text(input_path.as_path().to_str().unwrap())'
    '                        .size(16 * 2),'
    '                )'
    '                .numbering(NumberingId::new(2), IndentLevel::new(0))'
    '                .size(16 * 2),'
    '        );'
    '        let lines: Vec<String> = fs::read_to_string(input_path)'
    '            .unwrap()'
    '            .split("\n")'
    '            .map(str::to_string)'
    '            .collect();'
    '        for line in lines {'
)
$genFileBlockExtra = @(
    '            *self = self'
    '                .to_owned()'
    '                .add_paragraph(Paragraph::new().add_run(Run::new().add_text(line)));'
    '        }'
    ''
    '        Ok(())'
    '    }'
    '}'
    ''
    '// fn gen_file(input_path: PathBuf, doc: &mut Docx) -> Result<(), DocxError> {'
    '//     *doc = doc.to_owned().add_paragraph('
    '//         Paragraph::new()'
    '//             .add_run('
    '//                 Run::new()'
    '//                     .add_text(input_path.as_path().to_str().unwrap())'
    '//                     .size(16 * 2),'
    '//             )'
    '//             .numbering(NumberingId::new(2), IndentLevel::new(0))'
    '//             .size(16 * 2),'
    '//     );'
    '//     let lines: Vec<String> = fs::read_to_string(input_path)'
    '//         .unwrap()'
    '//         .split("\n")'
    '//         .map(str::to_string)'
    '//         .collect();'
    '//     for line in lines {'
    '//         *doc = doc'
    '//             .to_owned()'
    '//             .add_paragraph(Paragraph::new().add_run(Run::new().add_text(line)));'
    '//     }'
    ''
    '//     Ok(())'
    '// }'
)

$firstPara = 62
for ($i = 0; $i -lt $genFileBlock.Count; $i++) {
    $d.Paragraphs.Item($firstPara + $i).Range.Text = $genFileBlock[$i]
}

$insertAfter = $firstPara + $genFileBlock.Count - 1
foreach ($line in $genFileBlockExtra) {
    $d.Paragraphs.Item($insertAfter).Range.InsertParagraphAfter()
    $insertAfter = $insertAfter + 1
    $d.Paragraphs.Item($insertAfter).Range.Text = $line
}

# --- 4) Insert a blank (4-space) paragraph just above "    let mut doc = Docx::new();" ---
$d.Paragraphs.Item(43).Range.InsertParagraphBefore()
$d.Paragraphs.Item(43).Range.Text = '    '

Write-Output "edit complete"
